$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connectivity+Modem")

# --- Update existing row 2 message text to the generic signal-strength regex ---
$ws.Range("C2").Value = 'WAN:.*| signal .* (S.*) on port modem1: .*%, RSSI:-.*(dBm), SINR:.*(dB), RSRP:-.*(dB), RSRQ:-.*(dB), RFBAND: Band .*'

# --- Add the new "Problematic" column header (bold, vertical-top aligned) ---
$ws.Range("E1").Value = 'Problematic'
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").VerticalAlignment = -4160

# --- Mark existing rows 2-9 as Problematic = TRUE (vertical-top aligned) ---
foreach ($r in 2..9) {
    $ws.Range("E$r").Value = $true
    $ws.Range("E$r").VerticalAlignment = -4160
}

# --- Append the 5 new log-message rows ---
# Row 10
$ws.Range("A10").Value = 'WARN'
$ws.Range("B10").Value = 'WAN'
$ws.Range("C10").Value = 'WAN:.*|SIM error: NOSIM'
$ws.Range("D10").Value = 'The SIM card is not inserted, it is inserted improperly, either the SIM or the port for the SIM may be bad, or you have an incorrect size of SIM card. If the SIM is inserted properly, and is the correct size, and you are still getting the NO SIM error, then you need to test the SIM card in a different device to eliminate the SIM being the issue. https://customer.cradlepoint.com/s/article/What-type-of-SIM-cards-do-CradlePoint-modems-use '
$ws.Range("E10").Value = $true
$ws.Range("E10").VerticalAlignment = -4160
$ws.Rows.Item(10).RowHeight = 90

# Row 11
$ws.Range("A11").Value = 'WARN'
$ws.Range("B11").Value = 'WAN'
$ws.Range("C11").Value = 'suspending due to IP conflict '
$ws.Range("D11").Value = 'The Cradlepoint is detecting an IP Conflict, usually this happens when there is something on the network that has the same IP address as the one of the Cradlepoint''s networks.  The most common appearance of this is when the WAN connection is trying to use an IP Address that is within the same IP address of one of the Cradlepoint''s LANs.  To resolve that, change the subnet of your LAN, or get the WAN to give out addresses on a non-conflicting range. '
$ws.Range("E11").Value = $true
$ws.Range("E11").VerticalAlignment = -4160
$ws.Rows.Item(11).RowHeight = 90

# Row 12
$ws.Range("A12").Value = 'INFO'
$ws.Range("B12").Value = 'WAN'
$ws.Range("C12").Value = 'Connect Event: unknown error - state: connecting'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 'A CPPM unknown error usually indicates the modem or SIM is not functioning properly.  It can also indicate modem/SIM card provisioning errors.  '
$ws.Range("E12").Value = $true
$ws.Range("E12").VerticalAlignment = -4160
$ws.Rows.Item(12).RowHeight = 30

# Row 13
$ws.Range("A13").Value = 'INFO'
$ws.Range("B13").Value = 'WAN'
$ws.Range("C13").Value = 'state=connecting result=unknown error'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = 'A CPPM unknown error usually indicates the modem or SIM is not functioning properly.  It can also indicate modem/SIM card provisioning errors.  '
$ws.Range("E13").Value = $true
$ws.Range("E13").VerticalAlignment = -4160
$ws.Rows.Item(13).RowHeight = 30

# Row 14
$ws.Range("A14").Value = 'INFO'
$ws.Range("B14").Value = 'WAN'
$ws.Range("C14").Value = 'Module FW(.*) / SIM Carrier(.*) - mismatch'
$ws.Range("D14").Value = 'Indicates a mismatch between the modem firmware and the SIM. It is most likely to occur on modems that do not support Auto Carrier Selection. For modems that do not support Auto-Carrier Selection, ensure the correct firmware is loaded. For information about modem firmware management, click https://customer.cradlepoint.com/s/article/NCOS-Cradlepoint-Manual-Modem-Firmware-Update.'
$ws.Range("E14").Value = $true
$ws.Range("E14").VerticalAlignment = -4160
$ws.Rows.Item(14).RowHeight = 90

# --- Restore the frozen-pane view state to the new scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C10").Select()

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1
